$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.9179281773574478
$ws.Range("J2").Value = 0.9179281773574478
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.2189473333333334
$ws.Range("N2").Value = 0.656842
$ws.Range("O2").Value = 0.009402596261870986
$ws.Range("P2").Value = 0.009402596261870984
$ws.Range("Q2").Value = 0.1350418253762222
$ws.Range("R2").Value = 1.215376428386
$ws.Range("S2").Value = 0.008630908049087185
$ws.Range("T2").Value = 0.008630908049087185

$ws.Range("I3").Value = 0.9179281773574478
$ws.Range("J3").Value = 0.9179281773574478
$ws.Range("O3").Value = 0.8622887582286424
$ws.Range("P3").Value = 0.8622887582286423
$ws.Range("S3").Value = 0.7915191481966347
$ws.Range("T3").Value = 0.7915191481966346

$ws.Range("I4").Value = 0.9179281773574478
$ws.Range("J4").Value = 0.9179281773574478
$ws.Range("M4").Value = 2.823530666666667
$ws.Range("N4").Value = 8.470592
$ws.Range("O4").Value = 0.1212552739852724
$ws.Range("P4").Value = 0.1212552739852723
$ws.Range("Q4").Value = 1.741490656348444
$ws.Range("R4").Value = 15.673415907136
$ws.Range("S4").Value = 0.111303632644279
$ws.Range("T4").Value = 0.111303632644279

$ws.Range("I5").Value = 0.9179281773574478
$ws.Range("J5").Value = 0.9179281773574478
$ws.Range("M5").Value = 0.1642436666666667
$ws.Range("N5").Value = 0.492731
$ws.Range("O5").Value = 0.007053371524214274
$ws.Range("P5").Value = 0.007053371524214274
$ws.Range("Q5").Value = 0.1013018254914445
$ws.Range("R5").Value = 0.911716429423
$ws.Range("S5").Value = 0.006474488467446931
$ws.Range("T5").Value = 0.006474488467446931

$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.055146
$ws.Range("H6").Value = 0.165438
$ws.Range("I6").Value = 0.08207182264255215
$ws.Range("J6").Value = 0.08207182264255215
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.2189473333333334
$ws.Range("N6").Value = 0.656842
$ws.Range("O6").Value = 0.009402596261870986
$ws.Range("P6").Value = 0.009402596261870984
$ws.Range("Q6").Value = 0.012074069644
$ws.Range("R6").Value = 0.108666626796
$ws.Range("S6").Value = 0.0007716882127837995
$ws.Range("T6").Value = 0.0007716882127837992

$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.055146
$ws.Range("H7").Value = 0.165438
$ws.Range("I7").Value = 0.08207182264255215
$ws.Range("J7").Value = 0.08207182264255215
$ws.Range("O7").Value = 0.8622887582286424
$ws.Range("P7").Value = 0.8622887582286423
$ws.Range("Q7").Value = 1.107282949318
$ws.Range("R7").Value = 9.965546543862
$ws.Range("S7").Value = 0.07076961003200767
$ws.Range("T7").Value = 0.07076961003200766

$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.055146
$ws.Range("H8").Value = 0.165438
$ws.Range("I8").Value = 0.08207182264255215
$ws.Range("J8").Value = 0.08207182264255215
$ws.Range("M8").Value = 2.823530666666667
$ws.Range("N8").Value = 8.470592
$ws.Range("O8").Value = 0.1212552739852724
$ws.Range("P8").Value = 0.1212552739852723
$ws.Range("Q8").Value = 0.155706422144
$ws.Range("R8").Value = 1.401357799296
$ws.Range("S8").Value = 0.009951641340993342
$ws.Range("T8").Value = 0.00995164134099334

$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.055146
$ws.Range("H9").Value = 0.165438
$ws.Range("I9").Value = 0.08207182264255215
$ws.Range("J9").Value = 0.08207182264255215
$ws.Range("M9").Value = 0.1642436666666667
$ws.Range("N9").Value = 0.492731
$ws.Range("O9").Value = 0.007053371524214274
$ws.Range("P9").Value = 0.007053371524214274
$ws.Range("Q9").Value = 0.009057381242000001
$ws.Range("R9").Value = 0.081516431178
$ws.Range("S9").Value = 0.0005788830567673416
$ws.Range("T9").Value = 0.0005788830567673416
